$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("MAE") before the existing "Tipo" column,
# shifting "Tipo"/"single" to column E.
$ws.Range("D1").EntireColumn.Insert()

# Copy header formatting (bold, centered, bordered) from the neighboring
# "R2" header cell (C1) onto the new "MAE" header cell (D1).
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text.
$ws.Range("D1").Value = "MAE"

# Update data row 2: new MSE / R2 values, plus the new MAE value.
$ws.Range("B2").Value = 0.09583082560078589
$ws.Range("C2").Value = 0.9994291631146087
$ws.Range("D2").Value = 0.219217455348433

$excel.CutCopyMode = 0
